$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# The two sheets swap names, so rename through a temporary name to avoid a
# transient name clash ("hotel_info" <-> "review_info").
$ws1.Name = "__tmp_swap__"
$ws2.Name = "hotel_info"
$ws1.Name = "review_info"

# New header row for the review_info sheet (25 columns).
$reviewHeaders = @(
    "STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL",
    "Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date","review_title",
    "review_content","review_rating","trip_month","trip_purpose","value","rooms","Location",
    "Cleanliness","Sleep Quality","Service","Picture(yes=1)","respondent","response_date","response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# The old sheet had a second data row (hotel_info data) that does not belong
# on review_info any more - drop it, and drop any now-unused trailing columns.
$ws1.Range("A2:I2").ClearContents()

# New header row for the hotel_info sheet (10 columns, State added after Hotel_Name).
$hotelHeaders = @(
    "STR","Hotel_Name","State","City","Zip","TA_ReviewURL","Tripadvisor_Hotel_Name",
    "English_Reviews_num","Local_Rank","Total_Reviews_num"
)
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

# The single hotel data row, now with the Louisiana "State" value inserted.
$ws2.Cells.Item(2, 1).Value = 1010
$ws2.Cells.Item(2, 2).Value = "La Quinta Inns & Suites New Orleans Slidell"
$ws2.Cells.Item(2, 3).Value = "Louisiana"
$ws2.Cells.Item(2, 4).Value = "Slidell"
$ws2.Cells.Item(2, 5).Value = 70461
$ws2.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g40435-d89189-Reviews-La_Quinta_Inn_New_Orleans_Slidell-Slidell_Louisiana.html"
$ws2.Cells.Item(2, 7).Value = "La Quinta Inn New Orleans Slidell"
$ws2.Cells.Item(2, 8).Value = "545"
$ws2.Cells.Item(2, 9).Value = "11"
$ws2.Cells.Item(2, 10).Value = "552"

# Remove the now-unused trailing columns (J:Y) that used to hold the
# review_info header on this sheet.
$ws2.Range("J1:Y1").ClearContents()
